# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# price table on Sheet1 with freshly scraped values, matching the
# GitHub Actions cron job that refreshes this workbook.
#
# Values in columns D/E are free-form display text (e.g. "62.694.42",
# "  -0.85%  "), not numeric cells, even though some look numeric.
# A leading apostrophe forces Excel to store the assignment as literal
# text (quote-prefixed), exactly like typing it in by hand, so strings
# such as "1.00" or "0.998" are not silently coerced into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $ws.Range($range).Value = "'" + $text
}

# Row 2: Bitcoin
Set-TextValue "D2" "62.694.42"
Set-TextValue "E2" "  -0.85%  "

# Row 3: Ethereum
Set-TextValue "D3" "2.461.13"
Set-TextValue "E3" "  -0.14%  "

# Row 4: TetherUSD
Set-TextValue "E4" "  +0.10%  "

# Row 5: BNB
Set-TextValue "D5" "573.59"
Set-TextValue "E5" "  -0.98%  "

# Row 6: Solana
Set-TextValue "D6" "146.99"
Set-TextValue "E6" "  -0.12%  "

# Row 7: USDC
Set-TextValue "E7" "  +0.02%  "

# Row 8: XRP
Set-TextValue "E8" "  -1.73%  "

# Row 9: Dogecoin
Set-TextValue "E9" "  -0.92%  "

# Row 10: TRON
Set-TextValue "E10" "  -0.74%  "

# Row 11: Toncoin
Set-TextValue "E11" "  -0.90%  "

# Row 12: Cardano
Set-TextValue "E12" "  -0.59%  "

# Row 13: Avalanche
Set-TextValue "D13" "29.02"
Set-TextValue "E13" "  +1.87%  "

# Row 14: ShibaInu
Set-TextValue "E14" "  -2.80%  "

# Row 15: WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.908.55"
Set-TextValue "E15" "  +0.34%  "

# Row 16: WrappedBTC
Set-TextValue "D16" "62.616.29"
Set-TextValue "E16" "  -0.61%  "

# Row 17: WrappedEther
Set-TextValue "D17" "2.465.53"
Set-TextValue "E17" "  +0.28%  "

# Row 18: Uniswap
Set-TextValue "D18" "7.94"
Set-TextValue "E18" "  -1.02%  "

# Row 19: Chainlink
Set-TextValue "D19" "10.94"
Set-TextValue "E19" "  -1.57%  "

# Row 20: BitcoinCash
Set-TextValue "D20" "326.08"
Set-TextValue "E20" "  -1.58%  "

# Row 22: SuiNetwork
Set-TextValue "D22" "2.17"
Set-TextValue "E22" "  +0.74%  "

# Row 23: Dai
Set-TextValue "E23" "  -0.02%  "

# Row 24: Aptos
Set-TextValue "E24" "  +15.73%  "

# Row 25: Litecoin
Set-TextValue "D25" "65.35"
Set-TextValue "E25" "  -1.78%  "

# Row 26: Bittensor
Set-TextValue "D26" "638.26"
Set-TextValue "E26" "  -2.03%  "

# Row 27: WrappedeETH
Set-TextValue "D27" "2.581.30"
Set-TextValue "E27" "  +0.31%  "

# Row 28: PEPE
Set-TextValue "E28" "  -3.93%  "

# Row 29: Binance-PegBSC-USD
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  -20.90%  "

# Row 30: Fetch.AI
Set-TextValue "E30" "  -1.42%  "

# Row 31: InternetComputer(DFINITY)
Set-TextValue "D31" "7.93"
Set-TextValue "E31" "  -3.60%  "

# Row 32: PancakeSwap
Set-TextValue "E32" "  -2.87%  "

# Row 33: Kaspa
Set-TextValue "E33" "  -2.47%  "

# Row 34: FirstDigitalUSD
Set-TextValue "D34" "0.998"
Set-TextValue "E34" "  +0.03%  "

# Row 35: ImmutableX
Set-TextValue "E35" "  +1.77%  "

# Row 36: NEARProtocol
Set-TextValue "D36" "4.74"
Set-TextValue "E36" "  -1.55%  "

# Row 37: Monero
Set-TextValue "D37" "151.77"
Set-TextValue "E37" "  -0.54%  "

# Row 38: PolygonEcosystemToken
Set-TextValue "E38" "  -2.09%  "

# Row 39: EthereumClassic
Set-TextValue "D39" "18.61"
Set-TextValue "E39" "  -1.41%  "

# Row 40: RenderToken
Set-TextValue "D40" "5.32"
Set-TextValue "E40" "  -4.77%  "

# Row 41: dogwifhat
Set-TextValue "E41" "  -0.45%  "

# Row 42: Stacks
Set-TextValue "D42" "1.73"
Set-TextValue "E42" "  -3.00%  "

# Row 43: USDe
Set-TextValue "E43" "  -0.06%  "

# Row 44: BabyDogeCoin
Set-TextValue "D44" "0.0₆0302"
Set-TextValue "E44" "  -25.17%  "

# Row 45: Aave
Set-TextValue "D45" "152.82"
Set-TextValue "E45" "  +3.67%  "

# Row 46: WhiteBITCoin
Set-TextValue "D46" "15.25"
Set-TextValue "E46" "  +1.62%  "

# Row 47: Filecoin
Set-TextValue "D47" "3.57"
Set-TextValue "E47" "  -2.17%  "

# Row 48: Mantle
Set-TextValue "E48" "  -0.48%  "

# Row 49: InjectiveProtocol
Set-TextValue "D49" "20.28"
Set-TextValue "E49" "  -2.79%  "

# Row 50: Hedera
Set-TextValue "E50" "  -2.09%  "

# Row 51: Stellar
Set-TextValue "E51" "  -1.69%  "
